# Generate Report for Handback
# Updates the localization status workbook: the f3766040 file has been
# handed back (in sync with en-US), and the handback timestamps for the
# zh-cn / de-de target files are refreshed.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G2").Value = "2016-03-10 06:30:45"
$wsZhCn.Range("B3").Value = $status
$wsZhCn.Range("G3").Value = "2016-03-10 06:30:45"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("G2").Value = "2016-03-10 06:31:04"
$wsDeDe.Range("B3").Value = $status
$wsDeDe.Range("G3").Value = "2016-03-10 06:31:04"
